# "Commit with included Kramatorsk"
#
# The "Oblast Center" flag column (H) is replaced by its logical inverse,
# "Not Oblast Center": the header text changes and every 0/1 value in the
# column is flipped accordingly (cities that were flagged as an Oblast
# Center, e.g. Kramatorsk, are now correctly represented under the new
# "Not Oblast Center" semantics, and vice versa).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the column H header from "Oblast Center" to "Not Oblast Center".
$ws.Range("H1").Value = "Not Oblast Center"

# Column H holds a 0/1 flag for each of the 48 data rows (rows 2-49).
# Flip every value so the column now represents the inverse condition.
for ($r = 2; $r -le 49; $r++) {
    $cell = $ws.Cells.Item($r, 8)
    $current = $cell.Value()
    $cell.Value = 1 - $current
}

# Reflect the workbook's last on-screen selection.
$ws.Range("H5").Select()
